# Updated cryptos list on Sun Sep 10 05:02:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.107.93"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.642.88"
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "1.729.41"
$ws.Range("E12").Value = "  +5.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "26.107.48"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "189.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "144.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.130"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0482"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.876"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "1.123.38"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.519"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.06%  "
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.87"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.14"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0925"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("E51").Value = "  -1.06%  "
